# Release mCSD 3.9.0 with CP integrated
#
# Updates the "Metadata" sheet values (version bump, experimental flag,
# publish date, contact details, jurisdiction) and renames the second
# sheet tab, matching the upstream IG-publisher regeneration.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Version 3.8.0 -> 3.9.0
$ws1.Range("B3").Value = "3.9.0"

# Experimental flag is now explicitly populated ("false"). Writing the bare
# word through .Value would be auto-typed as a Boolean by the COM layer
# (matching genuine Excel semantics), so it is staged as text elsewhere
# (apostrophe-prefixed to force text entry) and copy/pasted in as a
# literal string to keep the cell a plain string cell (t="s"), exactly
# like the other metadata rows.
$scratch = $ws1.Range("Z1000")
$scratch.Value = "'false"
$scratch.Copy()
$ws1.Range("B7").PasteSpecial(-4163)
$scratch.Clear()

# Publish date refreshed for the 3.9.0 release
$ws1.Range("B8").Value = "2024-12-02T17:05:26-06:00"

# Contact block now carries three distinct rendered ContactDetail lines
# instead of the generic placeholder text
$ws1.Range("B10").Value = "null (https://www.ihe.net/ihe_domains/it_infrastructure/)"
$ws1.Range("B11").Value = "null (iti@ihe.net)"
$ws1.Range("B12").Value = "IHE IT Infrastructure Technical Committee (iti@ihe.net)"

# Jurisdiction now renders the full display text
$ws1.Range("B13").Value = "Global (Whole world)"

# Second tab is renamed by the IG publisher's generic "Include #N" naming
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Include #0"
